$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values for data rows 2-10; update per regenerated save_data.
$newValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
